$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "290.66"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "-3.63%"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "30.79"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "-6.03%"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "4.961"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "0.33%"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.07233"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "-6.85%"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.801"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "-9.13%"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "-2.17%"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "3.760"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "-1.06%"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.8963"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "-2.76%"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.1651"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "-6.66%"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07727"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "-1.96%"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.08024"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "-6.66%"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.03042"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "-3.34%"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "0.13%"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.001493"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "-1.51%"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.005703"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "-3.63%"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.467"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "0.15%"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.083"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "-3.30%"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.3313"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "-0.73%"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.1301"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "-1.29%"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "4.036"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "-5.83%"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.04508"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "-1.31%"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.001213"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "-1.02%"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "-9.35%"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.0001251"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "-0.06%"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01605"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "-6.14%"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.04406"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "-6.69%"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.007299"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "-7.57%"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1308"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "-3.41%"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.007687"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.001901"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "-18.51%"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "-12.74%"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.00005936"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "-5.14%"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.00000000750"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "-0.06%"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.247"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "173.92%"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "-3.25%"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.00002101"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "-0.06%"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0002001"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "-0.06%"
